$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Marking" row (row 11), "Right" column: 3 -> 5
$ws.Range("B11").Value = 5

# "Total" row (row 12), "Right" column: 24 -> 40
$ws.Range("B12").Value = 40

# "Total" row (row 12), "Max" column (Correct/Total marks): 19/84 -> 40/140
$ws.Range("E12").Value = "40/140"
